# Append the new match row (row 72) to Sheet1, mirroring the layout of the
# existing rows (row 71 is the last data row before this edit).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clone the formatting of the previous row (bordered/bold index cell style,
# datetime number format on the match-date column, etc.) onto the new row.
$ws.Range("A71:V71").Copy()
$ws.Range("A72:V72").PasteSpecial(-4122)  # xlPasteFormats

# Populate the new row's values.
$ws.Cells.Item(72, 1).Value  = 71
$ws.Cells.Item(72, 2).Value  = "azerbaijan"
$ws.Cells.Item(72, 3).Value  = "premier-league"
$ws.Cells.Item(72, 4).Value  = "2023-2024"
$ws.Cells.Item(72, 5).Value  = 45262.64583333334
$ws.Cells.Item(72, 6).Value  = "Turan"
$ws.Cells.Item(72, 7).Value  = 2
$ws.Cells.Item(72, 8).Value  = "Sumqayit"
$ws.Cells.Item(72, 9).Value  = 2
$ws.Cells.Item(72, 10).Value = 2.21
$ws.Cells.Item(72, 11).Value = "01/12/2023 03:43"
$ws.Cells.Item(72, 12).Value = 2.26
$ws.Cells.Item(72, 13).Value = "02/12/2023 15:26"
$ws.Cells.Item(72, 14).Value = 2.95
$ws.Cells.Item(72, 15).Value = "01/12/2023 03:43"
$ws.Cells.Item(72, 16).Value = 3.22
$ws.Cells.Item(72, 17).Value = "02/12/2023 15:26"
$ws.Cells.Item(72, 18).Value = 3.18
$ws.Cells.Item(72, 19).Value = "01/12/2023 03:43"
$ws.Cells.Item(72, 20).Value = 3.15
$ws.Cells.Item(72, 21).Value = "02/12/2023 15:26"
$ws.Cells.Item(72, 22).Value = "https://www.betexplorer.com/football/azerbaijan/premier-league/turan-sumqayit-fk/nN7PMxZl/"
